# Update the "Correct Answer" column (column 2) of the answer-key table
# for several questions, per the commit's corrected answer key.
#
# Table layout: row 1 is the header ("Question #" | "Correct Answer" | "Points"),
# and row N (N = question# + 1) holds the data for question #(N-1).
# Column 2 holds the "Correct Answer" value being changed here.

$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# Question 2: D -> B
$tbl.Cell(3, 2).Range.Text = "B"

# Question 3: True -> False
$tbl.Cell(4, 2).Range.Text = "False"

# Question 5: C -> D
$tbl.Cell(6, 2).Range.Text = "D"

# Question 7: B -> A
$tbl.Cell(8, 2).Range.Text = "A"

# Question 8: A -> C
$tbl.Cell(9, 2).Range.Text = "C"

# Question 10: False -> True
$tbl.Cell(11, 2).Range.Text = "True"

# Question 12: A -> D
$tbl.Cell(13, 2).Range.Text = "D"

# Question 13: D -> C
$tbl.Cell(14, 2).Range.Text = "C"

# Question 14: C -> A
$tbl.Cell(15, 2).Range.Text = "A"
